# Edit the SL fMRI experiment stimuli list:
#  - header row stays the same labels, but the underlying order in the shared
#    strings table changes as a side effect of the data edit (handled
#    automatically by just re-writing the header values)
#  - rows 2-25 get new "image"/"trialnum" values (stimuli were re-randomized
#    and trial numbers now continue a run, starting at 145 instead of 1)
#  - 24 new trial rows are appended (rows 26-49), extending the block to 192
#    trials total
#  - view/selection state is refreshed to match the re-saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("L.png",145,"R",0,1,1),
    @("B.png",146,"R",0,2,1),
    @("J.png",147,"R",0,3,1),
    @("B.png",148,"R",0,1,1),
    @("L.png",149,"R",0,2,1),
    @("A.png",150,"R",0,3,1),
    @("H.png",151,"R",0,1,1),
    @("M.png",152,"R",0,2,1),
    @("F.png",153,"R",0,3,1),
    @("E.png",154,"R",0,1,1),
    @("B.png",155,"R",0,2,1),
    @("C.png",156,"R",0,3,1),
    @("K.png",157,"R",0,1,1),
    @("A.png",158,"R",0,2,1),
    @("F.png",159,"R",0,3,1),
    @("M.png",160,"R",0,1,1),
    @("G.png",161,"R",0,2,1),
    @("E.png",162,"R",0,3,1),
    @("D.png",163,"R",0,1,1),
    @("M.png",164,"R",0,2,1),
    @("H.png",165,"R",0,3,1),
    @("K.png",166,"R",0,1,1),
    @("K.png",167,"R",0,2,1),
    @("J.png",168,"R",0,3,1),
    @("A.png",169,"R",0,1,1),
    @("H.png",170,"R",0,2,1),
    @("G.png",171,"R",0,3,1),
    @("C.png",172,"R",0,1,1),
    @("F.png",173,"R",0,2,1),
    @("K.png",174,"R",0,3,1),
    @("D.png",175,"R",0,1,1),
    @("C.png",176,"R",0,2,1),
    @("H.png",177,"R",0,3,1),
    @("E.png",178,"R",0,1,1),
    @("L.png",179,"R",0,2,1),
    @("A.png",180,"R",0,3,1),
    @("F.png",181,"R",0,1,1),
    @("J.png",182,"R",0,2,1),
    @("M.png",183,"R",0,3,1),
    @("E.png",184,"R",0,1,1),
    @("G.png",185,"R",0,2,1),
    @("C.png",186,"R",0,3,1),
    @("D.png",187,"R",0,1,1),
    @("D.png",188,"R",0,2,1),
    @("G.png",189,"R",0,3,1),
    @("L.png",190,"R",0,1,1),
    @("B.png",191,"R",0,2,1),
    @("J.png",192,"R",0,3,1)
)

# Header (unchanged text, just rewritten so the shared-strings table gets
# rebuilt the way a fresh save would)
$ws.Range("A1").Value = "image"
$ws.Range("B1").Value = "trialnum"
$ws.Range("C1").Value = "condition"
$ws.Range("D1").Value = "word"
$ws.Range("E1").Value = "location"
$ws.Range("F1").Value = "repetition"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Refresh view state: single-pane view (no frozen/scrolled topLeftCell),
# selection parked one row below the new data block.
$ws.Range("F52").Select()
